$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ID" column header (row 3) and the risk sequence numbers (rows 4-8,
# col B) are replaced with the "Risk" label, matching the text already
# used in B2/R2. Row 9 (risk #6) keeps its numeric value untouched.
$ws.Range("B3:B8").Value = "Risk"

# Refresh the view state left behind by the edit: zoomed out a bit, with
# the header block (B2:R9) selected.
$ws.Activate()
$excel.ActiveWindow.Zoom = 61
$ws.Range("B2:R9").Select()
